# Add session metrics for record removal
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Remove Incomplete Records")

$data = @{
    13 = @(71.39, 73.05, 74.47, 74.7, 71.16, 71.87, 75.41, 74.23, 75.650000000000006, 76.36)
    14 = @(75.27, 72.67, 71.37, 74.84, 72.23, 72.02, 72.67, 73.75, 72.02, 72.23)
    15 = @(71.599999999999994, 74.2, 69.8, 70, 69.400000000000006, 73, 73.2, 74.599999999999994, 74.400000000000006, 72.2)
    16 = @(73.42, 76.02, 73.61, 72.489999999999995, 70.45, 71.75, 72.3, 70.63, 71.930000000000007, 72.86)
    17 = @(77.08, 72.739999999999995, 72.22, 75, 72.22, 72.400000000000006, 73.44, 70.66, 73.44, 72.569999999999993)
    18 = @(73.5, 72.52, 74.63, 74.150000000000006, 72.2, 73.66, 75.61, 73.17, 74.47, 73.98)
    19 = @(75.040000000000006, 71.67, 69.37, 74.89, 73.05, 71.67, 72.59, 73.2, 73.81, 72.59)
    20 = @(73.7, 73.989999999999995, 75.430000000000007, 70.95, 73.55, 71.97, 74.13, 72.98, 71.53, 73.84)
    21 = @(72.19, 74.25, 73.56, 71.23, 72.47, 73.290000000000006, 73.84, 71.37, 70.41, 70.680000000000007)
}

foreach ($row in 13..21) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i  # C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$ws.Range("Q19").Select()
